$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "14-jun" column (G) right after the existing "13-jun" column (F),
# matching the formatting already used by the header/data in column F.

# Header cell G1 - copy F1's format (text-style header), then set the label.
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "14-jun"

# Data cells G2:G11 - copy F2:F11's format (centered integer), then fill values.
$ws.Range("F2:F11").Copy()
$ws.Range("G2:G11").PasteSpecial(-4122)

$ws.Range("G2").Value = 15
$ws.Range("G3").Value = 12
$ws.Range("G4").Value = 8
$ws.Range("G5").Value = 14
$ws.Range("G6").Value = 16
$ws.Range("G7").Value = 15
$ws.Range("G8").Value = 10
$ws.Range("G9").Value = 15
$ws.Range("G10").Value = 20
$ws.Range("G11").Value = 10

$excel.CutCopyMode = $false

$null = $ws.Range("G12").Select()
